# chore: update Sheets via scheduled runner
#
# Refreshes cached market-board figures (currentAveragePrice / NQ / HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the ALC,
# ARM, BSM, CRP, GSM, LTW and WVR leve-profit sheets, reflecting the latest
# scheduled data pull. No structural/formula changes; plain value updates.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2916.5386
$ws.Range("I100").Value = 2002.5
$ws.Range("J100").Value = 3700
$ws.Range("K100").Value = 2002.5
$ws.Range("L100").Value = 3700
$ws.Range("M100").Value = -1461.5
$ws.Range("N100").Value = -4782

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1095.9608
$ws.Range("I112").Value = 429.83334
$ws.Range("J112").Value = 1184.7778
$ws.Range("K112").Value = 1289.50002
$ws.Range("L112").Value = 3554.3334
$ws.Range("M112").Value = -181.5000199999999
$ws.Range("N112").Value = -5770.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 16668859
$ws.Range("I137").Value = 1466.6666
$ws.Range("J137").Value = 22224656
$ws.Range("K137").Value = 4399.9998
$ws.Range("L137").Value = 66673968
$ws.Range("M137").Value = -1849.9998
$ws.Range("N137").Value = -66679068

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3362.3044
$ws.Range("J138").Value = 3603.3333
$ws.Range("L138").Value = 10809.9999
$ws.Range("N138").Value = -21089.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12501798
$ws.Range("I61").Value = 13515268
$ws.Range("J61").Value = 2342.6667
$ws.Range("K61").Value = 13515268
$ws.Range("L61").Value = 2342.6667
$ws.Range("M61").Value = -13515056
$ws.Range("N61").Value = -2766.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 14288481
$ws.Range("I74").Value = 19232038
$ws.Range("J74").Value = 7095.1113
$ws.Range("K74").Value = 19232038
$ws.Range("L74").Value = 7095.1113
$ws.Range("M74").Value = -19231164
$ws.Range("N74").Value = -8843.1113

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 14288481
$ws.Range("I77").Value = 19232038
$ws.Range("J77").Value = 7095.1113
$ws.Range("K77").Value = 96160190
$ws.Range("L77").Value = 35475.5565
$ws.Range("M77").Value = -96155822
$ws.Range("N77").Value = -44211.5565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7144802.5
$ws.Range("I132").Value = 8622163
$ws.Range("K132").Value = 25866489
$ws.Range("M132").Value = -25863959

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 12501798
$ws.Range("I136").Value = 13515268
$ws.Range("J136").Value = 2342.6667
$ws.Range("K136").Value = 40545804
$ws.Range("L136").Value = 7028.000100000001
$ws.Range("M136").Value = -40543254
$ws.Range("N136").Value = -12128.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 813.9091
$ws.Range("I107").Value = 865.3
$ws.Range("K107").Value = 865.3
$ws.Range("M107").Value = 1054.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 60113.332
$ws.Range("J123").Value = 60113.332
$ws.Range("L123").Value = 60113.332
$ws.Range("N123").Value = -69913.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2428.4644
$ws.Range("I134").Value = 1360.3158
$ws.Range("J134").Value = 4683.4443
$ws.Range("K134").Value = 4080.9474
$ws.Range("L134").Value = 14050.3329
$ws.Range("M134").Value = -1545.9474
$ws.Range("N134").Value = -19120.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8337449.5
$ws.Range("I31").Value = 5477.64
$ws.Range("J31").Value = 22224070
$ws.Range("K31").Value = 5477.64
$ws.Range("L31").Value = 22224070
$ws.Range("M31").Value = -5182.64
$ws.Range("N31").Value = -22224660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8337449.5
$ws.Range("I34").Value = 5477.64
$ws.Range("J34").Value = 22224070
$ws.Range("K34").Value = 5477.64
$ws.Range("L34").Value = 22224070
$ws.Range("M34").Value = -5275.64
$ws.Range("N34").Value = -22224474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3363.5
$ws.Range("I58").Value = 623.2857
$ws.Range("J58").Value = 7199.8
$ws.Range("K58").Value = 623.2857
$ws.Range("L58").Value = 7199.8
$ws.Range("M58").Value = -420.2857
$ws.Range("N58").Value = -7605.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2547.8286
$ws.Range("I132").Value = 1815.931
$ws.Range("K132").Value = 5447.793
$ws.Range("M132").Value = -2917.793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1927.75
$ws.Range("I134").Value = 1798.381
$ws.Range("J134").Value = 2833.3333
$ws.Range("K134").Value = 5395.143
$ws.Range("L134").Value = 8499.999899999999
$ws.Range("M134").Value = -2860.143
$ws.Range("N134").Value = -13569.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3363.5
$ws.Range("I136").Value = 623.2857
$ws.Range("J136").Value = 7199.8
$ws.Range("K136").Value = 1869.8571
$ws.Range("L136").Value = 21599.4
$ws.Range("M136").Value = 680.1428999999998
$ws.Range("N136").Value = -26699.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 44479.266
$ws.Range("I70").Value = 153195
$ws.Range("J70").Value = 4946.273
$ws.Range("K70").Value = 153195
$ws.Range("L70").Value = 4946.273
$ws.Range("M70").Value = -152925
$ws.Range("N70").Value = -5486.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 44479.266
$ws.Range("I73").Value = 153195
$ws.Range("J73").Value = 4946.273
$ws.Range("K73").Value = 153195
$ws.Range("L73").Value = 4946.273
$ws.Range("M73").Value = -152259
$ws.Range("N73").Value = -6818.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16048410
$ws.Range("I80").Value = 30305758
$ws.Range("J80").Value = 2979173.5
$ws.Range("K80").Value = 30305758
$ws.Range("L80").Value = 2979173.5
$ws.Range("M80").Value = -30304760
$ws.Range("N80").Value = -2981169.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 16048410
$ws.Range("I83").Value = 30305758
$ws.Range("J83").Value = 2979173.5
$ws.Range("K83").Value = 151528790
$ws.Range("L83").Value = 14895867.5
$ws.Range("M83").Value = -151523798
$ws.Range("N83").Value = -14905851.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 67907.92999999999
$ws.Range("I113").Value = 143587.42
$ws.Range("J113").Value = 1688.375
$ws.Range("K113").Value = 143587.42
$ws.Range("L113").Value = 1688.375
$ws.Range("M113").Value = -141417.42
$ws.Range("N113").Value = -6028.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3989.442
$ws.Range("I132").Value = 3086.6667
$ws.Range("J132").Value = 5129.7896
$ws.Range("K132").Value = 9260.000100000001
$ws.Range("L132").Value = 15389.3688
$ws.Range("M132").Value = -6730.000100000001
$ws.Range("N132").Value = -20449.3688

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2788.375
$ws.Range("I82").Value = 2656.75
$ws.Range("J82").Value = 2920
$ws.Range("K82").Value = 2656.75
$ws.Range("L82").Value = 2920
$ws.Range("M82").Value = -2295.75
$ws.Range("N82").Value = -3642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2788.375
$ws.Range("I85").Value = 2656.75
$ws.Range("J85").Value = 2920
$ws.Range("K85").Value = 2656.75
$ws.Range("L85").Value = 2920
$ws.Range("M85").Value = -1408.75
$ws.Range("N85").Value = -5416

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11088.889
$ws.Range("I122").Value = 26950
$ws.Range("J122").Value = 6557.143
$ws.Range("K122").Value = 80850
$ws.Range("L122").Value = 19671.429
$ws.Range("M122").Value = -78400
$ws.Range("N122").Value = -24571.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12829459
$ws.Range("I132").Value = 5792.8335
$ws.Range("J132").Value = 33347326
$ws.Range("K132").Value = 17378.5005
$ws.Range("L132").Value = 100041978
$ws.Range("M132").Value = -14848.5005
$ws.Range("N132").Value = -100047038

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 19238244
$ws.Range("I136").Value = 31251272
$ws.Range("J136").Value = 17399.7
$ws.Range("K136").Value = 93753816
$ws.Range("L136").Value = 52199.10000000001
$ws.Range("M136").Value = -93751266
$ws.Range("N136").Value = -57299.10000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 50343.668
$ws.Range("J95").Value = 50343.668
$ws.Range("L95").Value = 50343.668
$ws.Range("N95").Value = -55835.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1226.9166
$ws.Range("I132").Value = 933.2
$ws.Range("J132").Value = 2108.0667
$ws.Range("K132").Value = 2799.6
$ws.Range("L132").Value = 6324.2001
$ws.Range("M132").Value = -269.6000000000004
$ws.Range("N132").Value = -11384.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1114
$ws.Range("I136").Value = 1023.62067
$ws.Range("J136").Value = 1550.8334
$ws.Range("K136").Value = 3070.86201
$ws.Range("L136").Value = 4652.5002
$ws.Range("M136").Value = -520.8620099999998
$ws.Range("N136").Value = -9752.5002

